$d = $word.ActiveDocument

$replacements = @(
    @{old="891×6="; new="881×8="},
    @{old="309×8="; new="989×5="},
    @{old="639×4="; new="273×2="},
    @{old="166×2="; new="894×3="},
    @{old="698×5="; new="721×5="},
    @{old="960×6="; new="754×2="},
    @{old="359×2="; new="285×9="},
    @{old="215×7="; new="720×6="},
    @{old="923×5="; new="359×7="},
    @{old="937×8="; new="790×5="},
    @{old="714×7="; new="988×2="},
    @{old="368×8="; new="131×4="},
    @{old="786×3="; new="525×4="},
    @{old="812×9="; new="897×6="},
    @{old="257×4="; new="219×6="},
    @{old="419×9="; new="532×3="},
    @{old="149×4="; new="194×4="},
    @{old="855×7="; new="216×3="},
    @{old="631×4="; new="757×3="},
    @{old="812×6="; new="925×7="},
    @{old="305×3="; new="869×8="},
    @{old="929×8="; new="930×6="},
    @{old="189×6="; new="388×7="},
    @{old="723×7="; new="274×5="},
    @{old="788×8="; new="220×4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying replacements"
